$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure D-column price cells that look numeric stay as plain text (matches source formatting)
$textCells = @("D2", "D3", "D5", "D6", "D7", "D9", "D10", "D11", "D13", "D14", "D15", "D16", "D17", "D18", "D19", "D20", "D21", "D22", "D23", "D24", "D26", "D28", "D29", "D30", "D31", "D32", "D33", "D34", "D36", "D37", "D38", "D39", "D40", "D41", "D42", "D43", "D44", "D45", "D46", "D47", "D48", "D49", "D50", "D51")
foreach ($addr in $textCells) {
    $ws.Range($addr).NumberFormat = "@"
}

$ws.Range("D2").Value = '43.705.76'
$ws.Range("E2").Value = '  -1.52%  '
$ws.Range("D3").Value = '2.182.68'
$ws.Range("E3").Value = '  -2.61%  '
$ws.Range("E4").Value = '  +0.18%  '
$ws.Range("D5").Value = '291.91'
$ws.Range("E5").Value = '  -4.74%  '
$ws.Range("D6").Value = '88.38'
$ws.Range("E6").Value = '  -4.66%  '
$ws.Range("D7").Value = '0.562'
$ws.Range("E7").Value = '  -1.47%  '
$ws.Range("E8").Value = '  +0.01%  '
$ws.Range("D9").Value = '0.474'
$ws.Range("E9").Value = '  -9.03%  '
$ws.Range("D10").Value = '31.97'
$ws.Range("E10").Value = '  -7.13%  '
$ws.Range("D11").Value = '0.0760'
$ws.Range("E11").Value = '  -6.03%  '
$ws.Range("E12").Value = '  -1.95%  '
$ws.Range("D13").Value = '6.66'
$ws.Range("E13").Value = '  -6.61%  '
$ws.Range("D14").Value = '2.527.93'
$ws.Range("E14").Value = '  -2.03%  '
$ws.Range("D15").Value = '2.271.89'
$ws.Range("E15").Value = '  -3.97%  '
$ws.Range("D16").Value = '12.85'
$ws.Range("E16").Value = '  -4.90%  '
$ws.Range("D17").Value = '0.755'
$ws.Range("E17").Value = '  -9.38%  '
$ws.Range("D18").Value = '43.471.42'
$ws.Range("E18").Value = '  -1.22%  '
$ws.Range("D19").Value = '0.0₃0866'
$ws.Range("E19").Value = '  -9.72%  '
$ws.Range("D20").Value = '5.76'
$ws.Range("E20").Value = '  -8.90%  '
$ws.Range("D21").Value = '10.61'
$ws.Range("E21").Value = '  -14.06%  '
$ws.Range("D22").Value = '62.29'
$ws.Range("E22").Value = '  -5.00%  '
$ws.Range("D23").Value = '226.39'
$ws.Range("E23").Value = '  -4.51%  '
$ws.Range("D24").Value = '2.74'
$ws.Range("E24").Value = '  -12.44%  '
$ws.Range("E25").Value = '  -0.06%  '
$ws.Range("D26").Value = '1.79'
$ws.Range("E26").Value = '  -9.01%  '
$ws.Range("E27").Value = '  -0.78%  '
$ws.Range("B28").Value = 'Cosmos'
$ws.Range("C28").Value = 'https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom'
$ws.Range("D28").Value = '9.04'
$ws.Range("E28").Value = '  -7.44%  '
$ws.Range("B29").Value = 'InjectiveProtocol'
$ws.Range("C29").Value = 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
$ws.Range("D29").Value = '34.85'
$ws.Range("E29").Value = '  -9.67%  '
$ws.Range("B30").Value = 'EthereumClassic'
$ws.Range("C30").Value = 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'
$ws.Range("D30").Value = '18.76'
$ws.Range("E30").Value = '  -6.06%  '
$ws.Range("B31").Value = 'Monero'
$ws.Range("C31").Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Range("D31").Value = '147.05'
$ws.Range("E31").Value = '  -4.07%  '
$ws.Range("D32").Value = '5.23'
$ws.Range("E32").Value = '  -11.75%  '
$ws.Range("D33").Value = '2.47'
$ws.Range("E33").Value = '  -6.04%  '
$ws.Range("D34").Value = '0.0718'
$ws.Range("E34").Value = '  -9.70%  '
$ws.Range("E35").Value = '  -3.48%  '
$ws.Range("D36").Value = '2.84'
$ws.Range("E36").Value = '  -8.46%  '
$ws.Range("D37").Value = '0.0990'
$ws.Range("E37").Value = '  -9.00%  '
$ws.Range("D38").Value = '1.62'
$ws.Range("E38").Value = '  -8.90%  '
$ws.Range("B39").Value = 'FirstDigitalUSD'
$ws.Range("C39").Value = 'https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd'
$ws.Range("D39").Value = '1.01'
$ws.Range("E39").Value = '  -0.17%  '
$ws.Range("B40").Value = 'NEARProtocol'
$ws.Range("C40").Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
$ws.Range("D40").Value = '3.03'
$ws.Range("E40").Value = '  -12.32%  '
$ws.Range("B41").Value = 'VeChain'
$ws.Range("C41").Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range("D41").Value = '0.0275'
$ws.Range("E41").Value = '  -8.28%  '
$ws.Range("B42").Value = 'RenderToken'
$ws.Range("C42").Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range("D42").Value = '3.43'
$ws.Range("E42").Value = '  -9.76%  '
$ws.Range("B43").Value = 'Celestia'
$ws.Range("C43").Value = 'https://coinranking.com/coin/YQcD0lBl7+celestia-tia'
$ws.Range("D43").Value = '12.89'
$ws.Range("E43").Value = '  -12.53%  '
$ws.Range("D44").Value = '1.727.57'
$ws.Range("E44").Value = '  -0.43%  '
$ws.Range("D45").Value = '1.61'
$ws.Range("E45").Value = '  -0.27%  '
$ws.Range("D46").Value = '67.07'
$ws.Range("E46").Value = '  -2.16%  '
$ws.Range("D47").Value = '72.46'
$ws.Range("E47").Value = '  -9.95%  '
$ws.Range("B48").Value = 'RocketPoolETH'
$ws.Range("C48").Value = 'https://coinranking.com/coin/QJZRUGyNI+rocketpooleth-reth'
$ws.Range("D48").Value = '2.411.21'
$ws.Range("E48").Value = '  -1.99%  '
$ws.Range("B49").Value = 'Algorand'
$ws.Range("C49").Value = 'https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo'
$ws.Range("D49").Value = '0.168'
$ws.Range("E49").Value = '  -11.75%  '
$ws.Range("D50").Value = '90.25'
$ws.Range("E50").Value = '  -9.03%  '
$ws.Range("D51").Value = '7.33'
$ws.Range("E51").Value = '  -9.81%  '
